$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fila 7 (Luis Rondón): "portero" (E7) acertado -> pasa a verde y suma 10 puntos
$ws.Range("E7").Interior.Color = 5287936   # verde (igual a las demas celdas acertadas)
$ws.Range("G7").Value = 10

# Fila 11 (German Rodriguez): "portero" (E11) ya no acertado -> pasa a rojo y resta los 10 puntos
$ws.Range("E11").Interior.Color = 255      # rojo (igual a las demas celdas no acertadas)
$ws.Range("G11").Value = 0

# Mover la seleccion activa a C14
[void]$ws.Range("C14").Select()
